# Applies the scheduled market-data refresh to the Durandal_Profits sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 281.57144
$ws.Range("I19").Value = 289
$ws.Range("J19").Value = 280.33334
$ws.Range("K19").Value = 289
$ws.Range("L19").Value = 280.33334
$ws.Range("M19").Value = -114
$ws.Range("N19").Value = -630.33334
# Row 113
$ws.Range("H113").Value = 1949.3334
$ws.Range("I113").Value = 1761.75
$ws.Range("J113").Value = 3450
$ws.Range("K113").Value = 1761.75
$ws.Range("L113").Value = 3450
$ws.Range("M113").Value = 1492.25
$ws.Range("N113").Value = -9958
# Row 129
$ws.Range("H129").Value = 984.3333
$ws.Range("I129").Value = 325.44446
$ws.Range("J129").Value = 1149.0555
$ws.Range("K129").Value = 976.33338
$ws.Range("L129").Value = 3447.1665
$ws.Range("M129").Value = 4023.66662
$ws.Range("N129").Value = -13447.1665
# Row 133
$ws.Range("H133").Value = 57278.57
$ws.Range("J133").Value = 57278.57
$ws.Range("L133").Value = 57278.57
$ws.Range("N133").Value = -67398.57000000001
# Row 135
$ws.Range("H135").Value = 1982.098
$ws.Range("I135").Value = 886.2954999999999
$ws.Range("J135").Value = 8870
$ws.Range("K135").Value = 7976.6595
$ws.Range("L135").Value = 79830
$ws.Range("M135").Value = -5441.6595
$ws.Range("N135").Value = -84900

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2895.3076
$ws.Range("I45").Value = 2875.7778
$ws.Range("J45").Value = 2939.25
$ws.Range("K45").Value = 2875.7778
$ws.Range("L45").Value = 2939.25
$ws.Range("M45").Value = -2498.7778
$ws.Range("N45").Value = -3693.25
# Row 61
$ws.Range("H61").Value = 1113.7709
$ws.Range("I61").Value = 712.6774
$ws.Range("K61").Value = 712.6774
$ws.Range("M61").Value = -500.6774
# Row 74
$ws.Range("H74").Value = 244824.97
$ws.Range("I74").Value = 278703.44
$ws.Range("J74").Value = 900
$ws.Range("K74").Value = 278703.44
$ws.Range("L74").Value = 900
$ws.Range("M74").Value = -277829.44
$ws.Range("N74").Value = -2648
# Row 77
$ws.Range("H77").Value = 244824.97
$ws.Range("I77").Value = 278703.44
$ws.Range("J77").Value = 900
$ws.Range("K77").Value = 1393517.2
$ws.Range("L77").Value = 4500
$ws.Range("M77").Value = -1389149.2
$ws.Range("N77").Value = -13236
# Row 88
$ws.Range("H88").Value = 4825
$ws.Range("I88").Value = 3146.5
$ws.Range("J88").Value = 5664.25
$ws.Range("K88").Value = 3146.5
$ws.Range("L88").Value = 5664.25
$ws.Range("M88").Value = -2740.5
$ws.Range("N88").Value = -6476.25
# Row 91
$ws.Range("H91").Value = 4825
$ws.Range("I91").Value = 3146.5
$ws.Range("J91").Value = 5664.25
$ws.Range("K91").Value = 3146.5
$ws.Range("L91").Value = 5664.25
$ws.Range("M91").Value = -1742.5
$ws.Range("N91").Value = -8472.25
# Row 122
$ws.Range("H122").Value = 26317838
$ws.Range("I122").Value = 34484916
$ws.Range("K122").Value = 103454748
$ws.Range("M122").Value = -103452298
# Row 132
$ws.Range("H132").Value = 1545.0938
$ws.Range("I132").Value = 1111.6923
$ws.Range("K132").Value = 3335.0769
$ws.Range("M132").Value = -805.0769
# Row 136
$ws.Range("H136").Value = 1113.7709
$ws.Range("I136").Value = 712.6774
$ws.Range("K136").Value = 2138.0322
$ws.Range("M136").Value = 411.9677999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 29795
$ws.Range("J53").Value = 29795
$ws.Range("L53").Value = 29795
$ws.Range("N53").Value = -30943
# Row 134
$ws.Range("H134").Value = 4748.184
$ws.Range("I134").Value = 783.03705
$ws.Range("J134").Value = 14480.818
$ws.Range("K134").Value = 2349.11115
$ws.Range("L134").Value = 43442.454
$ws.Range("M134").Value = 185.8888499999998
$ws.Range("N134").Value = -48512.454

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2005.2727
$ws.Range("I132").Value = 1757.5294
$ws.Range("J132").Value = 2847.6
$ws.Range("K132").Value = 5272.5882
$ws.Range("L132").Value = 8542.799999999999
$ws.Range("M132").Value = -2742.5882
$ws.Range("N132").Value = -13602.8
# Row 134
$ws.Range("H134").Value = 1167.081
$ws.Range("I134").Value = 1199.625
$ws.Range("J134").Value = 958.8
$ws.Range("K134").Value = 3598.875
$ws.Range("L134").Value = 2876.4
$ws.Range("M134").Value = -1063.875
$ws.Range("N134").Value = -7946.4

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1153.0555
$ws.Range("I5").Value = 904.2308
$ws.Range("J5").Value = 1800
$ws.Range("K5").Value = 2712.6924
$ws.Range("L5").Value = 5400
$ws.Range("M5").Value = -2600.6924
$ws.Range("N5").Value = -5624
# Row 109
$ws.Range("H109").Value = 1170
$ws.Range("I109").Value = 730
$ws.Range("K109").Value = 2190
$ws.Range("M109").Value = -1150
# Row 114
$ws.Range("H114").Value = 689.7143
$ws.Range("I114").Value = 471.33334
$ws.Range("J114").Value = 2000
$ws.Range("K114").Value = 1414.00002
$ws.Range("L114").Value = 6000
$ws.Range("M114").Value = 1839.99998
$ws.Range("N114").Value = -12508
# Row 135
$ws.Range("H135").Value = 1153.0555
$ws.Range("I135").Value = 904.2308
$ws.Range("J135").Value = 1800
$ws.Range("K135").Value = 8138.077200000001
$ws.Range("L135").Value = 16200
$ws.Range("M135").Value = -5603.077200000001
$ws.Range("N135").Value = -21270
# Row 139
$ws.Range("H139").Value = 1543.75
$ws.Range("I139").Value = 1335.7142
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 4007.1426
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1132.8574
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 34.714287
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 59
$ws.Range("K2").Value = 25
$ws.Range("L2").Value = 59
$ws.Range("M2").Value = 88
$ws.Range("N2").Value = -285
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
# Row 43
$ws.Range("H43").Value = 19009.5
$ws.Range("J43").Value = 19009.5
$ws.Range("L43").Value = 19009.5
$ws.Range("N43").Value = -19311.5
# Row 46
$ws.Range("H46").Value = 12042.625
$ws.Range("J46").Value = 19800
$ws.Range("L46").Value = 19800
$ws.Range("N46").Value = -20112
# Row 57
$ws.Range("H57").Value = 26000
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9180
# Row 70
$ws.Range("H70").Value = 18155886
$ws.Range("I70").Value = 37513640
$ws.Range("J70").Value = 7991.0625
$ws.Range("K70").Value = 37513640
$ws.Range("L70").Value = 7991.0625
$ws.Range("M70").Value = -37513370
$ws.Range("N70").Value = -8531.0625
# Row 73
$ws.Range("H73").Value = 18155886
$ws.Range("I73").Value = 37513640
$ws.Range("J73").Value = 7991.0625
$ws.Range("K73").Value = 37513640
$ws.Range("L73").Value = 7991.0625
$ws.Range("M73").Value = -37512704
$ws.Range("N73").Value = -9863.0625
# Row 80
$ws.Range("H80").Value = 2877.5293
$ws.Range("I80").Value = 2212.5
$ws.Range("J80").Value = 3468.6667
$ws.Range("K80").Value = 2212.5
$ws.Range("L80").Value = 3468.6667
$ws.Range("M80").Value = -1214.5
$ws.Range("N80").Value = -5464.6667
# Row 83
$ws.Range("H83").Value = 2877.5293
$ws.Range("I83").Value = 2212.5
$ws.Range("J83").Value = 3468.6667
$ws.Range("K83").Value = 11062.5
$ws.Range("L83").Value = 17343.3335
$ws.Range("M83").Value = -6070.5
$ws.Range("N83").Value = -27327.3335
# Row 126
$ws.Range("H126").Value = 23810824
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 27779078
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 83337234
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -83342174
# Row 132
$ws.Range("H132").Value = 1755.0613
$ws.Range("I132").Value = 1720.2927
$ws.Range("J132").Value = 1933.25
$ws.Range("K132").Value = 5160.8781
$ws.Range("L132").Value = 5799.75
$ws.Range("M132").Value = -2630.8781
$ws.Range("N132").Value = -10859.75
# Row 141
$ws.Range("H141").Value = 43576.332
$ws.Range("J141").Value = 43576.332
$ws.Range("L141").Value = 43576.332
$ws.Range("N141").Value = -53936.332

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 603.19446
$ws.Range("I55").Value = 631.7083
$ws.Range("J55").Value = 546.1667
$ws.Range("K55").Value = 631.7083
$ws.Range("L55").Value = 546.1667
$ws.Range("M55").Value = -458.7083
$ws.Range("N55").Value = -892.1667
# Row 68
$ws.Range("H68").Value = 1826.2106
$ws.Range("I68").Value = 1748
$ws.Range("J68").Value = 1933.75
$ws.Range("K68").Value = 1748
$ws.Range("L68").Value = 1933.75
$ws.Range("M68").Value = -999
$ws.Range("N68").Value = -3431.75
# Row 71
$ws.Range("H71").Value = 1826.2106
$ws.Range("I71").Value = 1748
$ws.Range("J71").Value = 1933.75
$ws.Range("K71").Value = 8740
$ws.Range("L71").Value = 9668.75
$ws.Range("M71").Value = -4996
$ws.Range("N71").Value = -17156.75
# Row 132
$ws.Range("H132").Value = 8398.799999999999
$ws.Range("I132").Value = 9801
$ws.Range("J132").Value = 2790
$ws.Range("K132").Value = 29403
$ws.Range("L132").Value = 8370
$ws.Range("M132").Value = -26873
$ws.Range("N132").Value = -13430
# Row 133
$ws.Range("H133").Value = 75914.82000000001
$ws.Range("J133").Value = 75914.82000000001
$ws.Range("L133").Value = 75914.82000000001
$ws.Range("N133").Value = -80974.82000000001
# Row 136
$ws.Range("H136").Value = 3206.1904
$ws.Range("I136").Value = 2575.1924
$ws.Range("J136").Value = 4231.5625
$ws.Range("K136").Value = 7725.5772
$ws.Range("L136").Value = 12694.6875
$ws.Range("M136").Value = -5175.5772
$ws.Range("N136").Value = -17794.6875

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 25035000
$ws.Range("I2").Value = 50000000
$ws.Range("K2").Value = 50000000
$ws.Range("M2").Value = -49999888
# Row 107
$ws.Range("H107").Value = 585.6923
$ws.Range("I107").Value = 636.2222
$ws.Range("K107").Value = 1908.6666
$ws.Range("M107").Value = 11.33339999999998
# Row 126
$ws.Range("H126").Value = 45456120
$ws.Range("I126").Value = 58824680
$ws.Range("J126").Value = 3001
$ws.Range("K126").Value = 176474040
$ws.Range("L126").Value = 9003
$ws.Range("M126").Value = -176471570
$ws.Range("N126").Value = -13943
# Row 132
$ws.Range("H132").Value = 18383306
$ws.Range("I132").Value = 22728156
$ws.Range("J132").Value = 1251.9231
$ws.Range("K132").Value = 68184468
$ws.Range("L132").Value = 3755.7693
$ws.Range("M132").Value = -68181938
$ws.Range("N132").Value = -8815.7693
# Row 135
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140
# Row 136
$ws.Range("H136").Value = 668.4211
$ws.Range("I136").Value = 536.25
$ws.Range("J136").Value = 815.2778
$ws.Range("K136").Value = 1608.75
$ws.Range("L136").Value = 2445.8334
$ws.Range("M136").Value = 941.25
$ws.Range("N136").Value = -7545.8334

Write-Output "Updated 284 cells across 8 sheets."
